$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Start clean: wipe all existing cell content/formatting on the sheet ---
$ws.Cells.Clear()

# --- Header row (bold, no fill) ---
$ws.Range("A1").Value = "Number of Cones "
$ws.Range("B1").Value = "No bands"
$ws.Range("C1").Value = "Two bands"
$ws.Range("D1").Value = "Four Bands"
$ws.Range("E1").Value = "Six bands"
$ws.Range("A1:E1").Font.Bold = $true

# --- Column A data rows: 0 .. 16 ---
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15
$ws.Range("A18").Value = 16

# --- Single "Yes" answer cell ---
$ws.Range("D2").Value = "Yes"

# --- Column widths (character units, COM ColumnWidth) ---
# Column A width is unchanged from the original workbook, so it is left alone.
$ws.Columns.Item(2).ColumnWidth = 10.002
$ws.Columns.Item(3).ColumnWidth = 11.836
$ws.Columns.Item(4).ColumnWidth = 10.168
$ws.Columns.Item(5).ColumnWidth = 9.502

# --- Selection ---
$null = $ws.Range("D3").Select()
